$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Range("B2").Value = 27.81570239103635
$ws.Range("C2").Value = 30.80390501693298
$ws.Range("B3").Value = 27.82998793290199
$ws.Range("C3").Value = 30.82642585885509
$ws.Range("B4").Value = 27.81040101797698
$ws.Range("C4").Value = 30.80843604591822
$ws.Range("B5").Value = 27.83059714135541
$ws.Range("C5").Value = 30.82348619399256
$ws.Range("B6").Value = 27.81276892523978
$ws.Range("C6").Value = 30.81269365560632
$ws.Range("B7").Value = 27.77995967987434
$ws.Range("C7").Value = 30.8159316451517
$ws.Range("B8").Value = 27.81054231798197
$ws.Range("C8").Value = 30.81781439205926
$ws.Range("B9").Value = 27.79302848503535
$ws.Range("C9").Value = 30.8269700018248
$ws.Range("B10").Value = 27.77886993060811
$ws.Range("C10").Value = 30.82344178040631
$ws.Range("B11").Value = 27.79238803733759
$ws.Range("C11").Value = 30.81013083384278
$ws.Range("B12").Value = 27.80542458593479
$ws.Range("C12").Value = 30.816923542459

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Range("B2").Value = 19.71169903920278
$ws.Range("C2").Value = 27.18949078542111
$ws.Range("B3").Value = 19.72321208847485
$ws.Range("C3").Value = 27.14313349217949
$ws.Range("B4").Value = 19.70485216020915
$ws.Range("C4").Value = 27.19447114028931
$ws.Range("B5").Value = 19.69738454128513
$ws.Range("C5").Value = 27.21594110718867
$ws.Range("B6").Value = 19.71785702247167
$ws.Range("C6").Value = 27.18830513438048
$ws.Range("B7").Value = 19.70378262439458
$ws.Range("C7").Value = 27.18923036513203
$ws.Range("B8").Value = 19.717695462581
$ws.Range("C8").Value = 27.20481006915502
$ws.Range("B9").Value = 19.70897206965194
$ws.Range("C9").Value = 27.18783770036239
$ws.Range("B10").Value = 19.70175579881327
$ws.Range("C10").Value = 27.15716340365552
$ws.Range("B11").Value = 19.70801303412406
$ws.Range("C11").Value = 27.1798500604582
$ws.Range("B12").Value = 19.70952238412085
$ws.Range("C12").Value = 27.18502332582222

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "sigma_050"

$newSheet.Range("A1").Value = "Rows"
$newSheet.Range("B1").Value = "Noisy"
$newSheet.Range("C1").Value = "NLM-LBP"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = 14.79279422790961
$newSheet.Range("C2").Value = 22.96662184382448
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = 14.79902987844351
$newSheet.Range("C3").Value = 23.00401240077766
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = 14.79181344603722
$newSheet.Range("C4").Value = 22.99301473913094
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = 14.77823956324047
$newSheet.Range("C5").Value = 22.98628552826259
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = 14.80343918904243
$newSheet.Range("C6").Value = 22.99837574781676
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = 14.78719751716943
$newSheet.Range("C7").Value = 22.97308713479629
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = 14.80702704989456
$newSheet.Range("C8").Value = 22.99628064475888
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = 14.79418528960567
$newSheet.Range("C9").Value = 23.00637875222865
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = 14.7941084626089
$newSheet.Range("C10").Value = 22.97217952052352
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = 14.78041386147308
$newSheet.Range("C11").Value = 22.98003313687737
$newSheet.Range("A12").Value = "Média"
$newSheet.Range("B12").Value = 14.79282484854249
$newSheet.Range("C12").Value = 22.98762694489972
